$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update individual cell values in existing rows (2, 3, 4) ---

# Row 2
$ws.Range("K2").Value = 2.1
$ws.Range("U2").Value = 1.91
$ws.Range("V2").Value = 1.8
$ws.Range("W2").Value = 11
$ws.Range("AC2").Value = 8.5
$ws.Range("AM2").Value = 29
$ws.Range("AR2").Value = 101
$ws.Range("AS2").Value = 251
$ws.Range("AX2").Value = 10
$ws.Range("AZ2").Value = 34

# Row 3
$ws.Range("G3").Value = 1.87
$ws.Range("H3").Value = 3.05
$ws.Range("I3").Value = 4.55
$ws.Range("J3").Value = 2.45
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 4.9
$ws.Range("M3").Value = 1.1
$ws.Range("Q3").Value = 2.25
$ws.Range("R3").Value = 1.6
$ws.Range("S3").Value = 1.47
$ws.Range("T3").Value = 2.52
$ws.Range("U3").Value = 1.95
$ws.Range("X3").Value = 8
$ws.Range("Z3").Value = 16
$ws.Range("AA3").Value = 16.5
$ws.Range("AD3").Value = 6
$ws.Range("AE3").Value = 16.5
$ws.Range("AF3").Value = 90
$ws.Range("AJ3").Value = 15
$ws.Range("AL3").Value = 50
$ws.Range("AM3").Value = 60
$ws.Range("AN3").Value = 3.65
$ws.Range("AO3").Value = 9.5
$ws.Range("AP3").Value = 19
$ws.Range("AQ3").Value = 35
$ws.Range("AR3").Value = 70
$ws.Range("AT3").Value = 2.52
$ws.Range("AU3").Value = 7.3
$ws.Range("AV3").Value = 70
$ws.Range("AW3").Value = 6.2
$ws.Range("AX3").Value = 28
$ws.Range("AY3").Value = 32
$ws.Range("AZ3").Value = 175
$ws.Range("BA3").Value = 200
$ws.Range("BB3").Value = 450

# Row 4
$ws.Range("H4").Value = 2.65
$ws.Range("I4").Value = 2.62
$ws.Range("J4").Value = 3.75
$ws.Range("P4").Value = 2.4
$ws.Range("Q4").Value = 2.5
$ws.Range("R4").Value = 1.47
$ws.Range("U4").Value = 1.95
$ws.Range("V4").Value = 1.75
$ws.Range("W4").Value = 7.3
$ws.Range("AE4").Value = 15
$ws.Range("AG4").Value = 900
$ws.Range("AH4").Value = 6.5
$ws.Range("AI4").Value = 12
$ws.Range("AO4").Value = 18.5
$ws.Range("AP4").Value = 27
$ws.Range("AS4").Value = 400
$ws.Range("AU4").Value = 6.9
$ws.Range("AW4").Value = 4.4
$ws.Range("AZ4").Value = 70

# --- Delete row 7 (Al Qadisiya vs Al Feiha), shifting rows 8-9 up ---
$ws.Rows("7").Delete()
